$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for column A (Id), Q (Ost), R (Nord) across rows 25-32 and 35-36.
$updates = @{
    25 = @{ A = 112183150; Q = 572358; R = 6714972 }
    26 = @{ A = 112183134; Q = 572354; R = 6714968 }
    27 = @{ A = 112183151; Q = 572361; R = 6714978 }
    28 = @{ A = 112183148; Q = 572357; R = 6714903 }
    29 = @{ A = 112183143; Q = 572359; R = 6714905 }
    30 = @{ A = 112183145; Q = 572351; R = 6714907 }
    31 = @{ A = 112183146; Q = 572346; R = 6714917 }
    32 = @{ A = 112183141; Q = 572361; R = 6714980 }
    35 = @{ A = 112183140; Q = 572350; R = 6714962 }
    36 = @{ A = 112183149; Q = 572345; R = 6714965 }
}

foreach ($rowNum in $updates.Keys) {
    $vals = $updates[$rowNum]
    $ws.Range("A$rowNum").Value = $vals.A
    $ws.Range("Q$rowNum").Value = $vals.Q
    $ws.Range("R$rowNum").Value = $vals.R
}
